$wb = $excel.ActiveWorkbook

# --- Fill in bus geodata (x/y coordinates for each bus) ---
$wsGeo = $wb.Worksheets.Item("bus_geodata")

$geoData = @(
    @(0, 1, 0),
    @(1, 1, 1),
    @(2, 0, 2),
    @(3, 0, 3),
    @(4, 0, 4),
    @(5, 0, 5),
    @(6, 0, 6),
    @(7, 0, 7),
    @(8, 2, 2),
    @(9, 2, 3)
)

for ($i = 0; $i -lt $geoData.Length; $i++) {
    $r = $i + 2
    $row = $geoData[$i]
    $wsGeo.Range("A$r").Value = $row[0]
    $wsGeo.Range("B$r").Value = $row[1]
    $wsGeo.Range("C$r").Value = $row[2]
}

# carry the index-column styling (bold, bordered) down through the new rows,
# and give column E the same blank styled look
$wsGeo.Range("A2").Copy()
$wsGeo.Range("A3:A11").PasteSpecial(-4122)
$wsGeo.Range("E2:E11").PasteSpecial(-4122)

# --- Update selections left on other sheets while reviewing the data ---
$wsBus = $wb.Worksheets.Item("bus")
$wsBus.Range("A2:B11").Select()

$wsDtypes = $wb.Worksheets.Item("dtypes")
$wsDtypes.Range("I221").Select()

# --- Finish on the bus_geodata sheet, which becomes the active tab ---
$wsGeo.Activate()
$wsGeo.Range("C12").Select()
